$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "'246.02"
$ws.Range("D4").Value = "'5.365"
$ws.Range("D5").Value = "'0.05876"
$ws.Range("D6").Value = "'3.389"
$ws.Range("D7").Value = "'6.387"
$ws.Range("D8").Value = "'0.8137"
$ws.Range("D9").Value = "'0.9578"
$ws.Range("D10").Value = "'0.1420"
$ws.Range("D11").Value = "'0.03565"
$ws.Range("D12").Value = "'0.07322"
$ws.Range("D13").Value = "'0.03035"
$ws.Range("D14").Value = "'4.414"
$ws.Range("D15").Value = "'0.09391"
$ws.Range("D17").Value = "'0.04828"
$ws.Range("E18").Value = "'17OneONE"
$ws.Range("D19").Value = "'0.006108"
$ws.Range("D21").Value = "'0.0009812"
$ws.Range("D22").Value = "'0.00009699"
$ws.Range("D24").Value = "'2.200"
$ws.Range("D26").Value = "'0.1287"
$ws.Range("D27").Value = "'0.0002471"
$ws.Range("D40").Value = "'0.03858"
$ws.Range("D41").Value = "'0.006619"
$ws.Range("D42").Value = "'0.1074"
$ws.Range("D44").Value = "'0.005773"
$ws.Range("D45").Value = "'0.00005646"
$ws.Range("D48").Value = "'0.03448"
$ws.Range("E48").Value = "'47BOLOBOLOWorstin24h"
